$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '303.97'
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '0.04%'
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '35.54'
$ws.Cells.Item(3, 4).Style = "Normal"

$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '-4.47%'
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '5.053'
$ws.Cells.Item(4, 4).Style = "Normal"

$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '0.17%'
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.07863'
$ws.Cells.Item(5, 4).Style = "Normal"

$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '0.16%'
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '2.133'
$ws.Cells.Item(6, 4).Style = "Normal"

$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '-3.27%'
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '-0.64%'
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 2).Value = 'MXToken'

$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.9239'
$ws.Cells.Item(8, 4).Style = "Normal"

$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '-0.28%'
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 2).Value = 'LiechtensteinCryptoassetsExchange'

$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.09699'
$ws.Cells.Item(9, 4).Style = "Normal"

$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '-1.37%'
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 2).Value = 'WazirX'

$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.1845'
$ws.Cells.Item(10, 4).Style = "Normal"

$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '-1.94%'
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'

$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.08695'
$ws.Cells.Item(11, 4).Style = "Normal"

$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '0.36%'
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 2).Value = 'BitrueCoin'

$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.03581'
$ws.Cells.Item(12, 4).Style = "Normal"

$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '-2.33%'
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 2).Value = 'BitMartToken'

$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.09929'
$ws.Cells.Item(13, 4).Style = "Normal"

$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '-0.29%'
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 2).Value = 'BitForexToken'

$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.001433'
$ws.Cells.Item(14, 4).Style = "Normal"

$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '-4.68%'
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 2).Value = 'TigerCash'

$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.005666'
$ws.Cells.Item(15, 4).Style = "Normal"

$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '-0.42%'
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 2).Value = 'LEO'

$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.477'
$ws.Cells.Item(16, 4).Style = "Normal"

$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '0.42%'
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 2).Value = 'GateToken'

$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '4.136'
$ws.Cells.Item(17, 4).Style = "Normal"

$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '2.88%'
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 2).Value = 'BTSEToken'

$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.753'
$ws.Cells.Item(18, 4).Style = "Normal"

$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '22.17%'
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 2).Value = 'BitpandaEcosystemToken'

$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.3372'
$ws.Cells.Item(19, 4).Style = "Normal"

$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '-1.13%'
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 2).Value = 'ProBitToken'

$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.1348'
$ws.Cells.Item(20, 4).Style = "Normal"

$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '3.38%'
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 2).Value = 'MCDex'

$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.172'
$ws.Cells.Item(21, 4).Style = "Normal"

$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '8.51%'
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 2).Value = 'ZBToken'

$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.2210'
$ws.Cells.Item(22, 4).Style = "Normal"

$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '0.13%'
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 2).Value = 'CoinExToken'

$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.04572'
$ws.Cells.Item(23, 4).Style = "Normal"

$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '-0.66%'
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '-1.98%'
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.004833'
$ws.Cells.Item(25, 4).Style = "Normal"

$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '8.00%'
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.0001303'
$ws.Cells.Item(26, 4).Style = "Normal"

$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '-7.12%'
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '74.40%'
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01852'
$ws.Cells.Item(39, 4).Style = "Normal"

$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '0.43%'
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.04723'
$ws.Cells.Item(40, 4).Style = "Normal"

$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '-0.75%'
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.007795'
$ws.Cells.Item(41, 4).Style = "Normal"

$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '-3.08%'
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '-0.78%'
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.007763'
$ws.Cells.Item(43, 4).Style = "Normal"

$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '2.54%'
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.002226'
$ws.Cells.Item(44, 4).Style = "Normal"

$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '5.88%'
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.01134'
$ws.Cells.Item(45, 4).Style = "Normal"

$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '8.94%'
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.00006297'
$ws.Cells.Item(46, 4).Style = "Normal"

$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '0.22%'
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '-0.09%'
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '-0.11%'
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '50.61'
$ws.Cells.Item(49, 4).Style = "Normal"

$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '32.45%'
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.001903'
$ws.Cells.Item(50, 4).Style = "Normal"

$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '-29.50%'
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '-0.09%'
$ws.Cells.Item(51, 5).Style = "Normal"
